$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-13 down to 7-14
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly data point
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 45044
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107011
$ws.Range("J6").Value = "Tuna"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("Q6").Value = "$/caja 18 kilos"
$ws.Range("R6").Value = "Provincia de Melipilla"
$ws.Range("S6").Value = 972
$ws.Range("T6").Value = 18

Write-Output "done"
